$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column (D) to remain stored as text, matching the
# original inline-string cell type, while keeping the default (unstyled) look.

$cell = $ws.Cells.Item(2, 4)
$cell.NumberFormat = "@"
$cell.Value = '26.589.04'
$cell.Style = "Normal"
$ws.Cells.Item(2, 5).Value = '  -7.47%  '

$cell = $ws.Cells.Item(3, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.687.94'
$cell.Style = "Normal"
$ws.Cells.Item(3, 5).Value = '  -6.47%  '

$cell = $ws.Cells.Item(4, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.005'
$cell.Style = "Normal"
$ws.Cells.Item(4, 5).Value = '  +0.22%  '

$cell = $ws.Cells.Item(5, 4)
$cell.NumberFormat = "@"
$cell.Value = '216.44'
$cell.Style = "Normal"
$ws.Cells.Item(5, 5).Value = '  -6.54%  '

$cell = $ws.Cells.Item(6, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.005'
$cell.Style = "Normal"
$ws.Cells.Item(6, 5).Value = '  +0.12%  '

$cell = $ws.Cells.Item(7, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.4969'
$cell.Style = "Normal"
$ws.Cells.Item(7, 5).Value = '  -16.44%  '

$cell = $ws.Cells.Item(8, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.2606'
$cell.Style = "Normal"
$ws.Cells.Item(8, 5).Value = '  -6.34%  '

$cell = $ws.Cells.Item(9, 4)
$cell.NumberFormat = "@"
$cell.Value = '21.63'
$cell.Style = "Normal"
$ws.Cells.Item(9, 5).Value = '  -7.81%  '

$cell = $ws.Cells.Item(10, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.06103'
$cell.Style = "Normal"
$ws.Cells.Item(10, 5).Value = '  -11.07%  '

$cell = $ws.Cells.Item(11, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.07262'
$cell.Style = "Normal"
$ws.Cells.Item(11, 5).Value = '  -3.84%  '

$cell = $ws.Cells.Item(12, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.709.74'
$cell.Style = "Normal"
$ws.Cells.Item(12, 5).Value = '  -5.39%  '

$cell = $ws.Cells.Item(13, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.433'
$cell.Style = "Normal"
$ws.Cells.Item(13, 5).Value = '  -6.26%  '

$cell = $ws.Cells.Item(14, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.5713'
$cell.Style = "Normal"
$ws.Cells.Item(14, 5).Value = '  -9.09%  '

$cell = $ws.Cells.Item(15, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.917.60'
$cell.Style = "Normal"
$ws.Cells.Item(15, 5).Value = '  -6.47%  '

$cell = $ws.Cells.Item(16, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.000008239'
$cell.Style = "Normal"
$ws.Cells.Item(16, 5).Value = '  -11.40%  '

$cell = $ws.Cells.Item(17, 4)
$cell.NumberFormat = "@"
$cell.Value = '64.57'
$cell.Style = "Normal"
$ws.Cells.Item(17, 5).Value = '  -14.48%  '

$cell = $ws.Cells.Item(18, 4)
$cell.NumberFormat = "@"
$cell.Value = '26.634.85'
$cell.Style = "Normal"
$ws.Cells.Item(18, 5).Value = '  -7.20%  '

$cell = $ws.Cells.Item(19, 4)
$cell.NumberFormat = "@"
$cell.Value = '4.995'
$cell.Style = "Normal"
$ws.Cells.Item(19, 5).Value = '  -9.03%  '

$cell = $ws.Cells.Item(20, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.005'
$cell.Style = "Normal"
$ws.Cells.Item(20, 5).Value = '  +0.14%  '

$cell = $ws.Cells.Item(22, 4)
$cell.NumberFormat = "@"
$cell.Value = '182.14'
$cell.Style = "Normal"
$ws.Cells.Item(22, 5).Value = '  -13.93%  '

$cell = $ws.Cells.Item(23, 4)
$cell.NumberFormat = "@"
$cell.Value = '6.167'
$cell.Style = "Normal"
$ws.Cells.Item(23, 5).Value = '  -10.32%  '

$cell = $ws.Cells.Item(24, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.005'
$cell.Style = "Normal"
$ws.Cells.Item(24, 5).Value = '  +0.13%  '

$cell = $ws.Cells.Item(25, 4)
$cell.NumberFormat = "@"
$cell.Value = '144.81'
$cell.Style = "Normal"
$ws.Cells.Item(25, 5).Value = '  -6.27%  '

$cell = $ws.Cells.Item(26, 4)
$cell.NumberFormat = "@"
$cell.Value = '7.552'
$cell.Style = "Normal"
$ws.Cells.Item(26, 5).Value = '  -3.94%  '

$cell = $ws.Cells.Item(27, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.1130'
$cell.Style = "Normal"
$ws.Cells.Item(27, 5).Value = '  -11.49%  '

$cell = $ws.Cells.Item(28, 4)
$cell.NumberFormat = "@"
$cell.Value = '15.27'
$cell.Style = "Normal"
$ws.Cells.Item(28, 5).Value = '  -7.14%  '

$cell = $ws.Cells.Item(29, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.314'
$cell.Style = "Normal"
$ws.Cells.Item(29, 5).Value = '  -9.06%  '

$cell = $ws.Cells.Item(30, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.05591'
$cell.Style = "Normal"
$ws.Cells.Item(30, 5).Value = '  -9.83%  '

$cell = $ws.Cells.Item(31, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.319'
$cell.Style = "Normal"
$ws.Cells.Item(31, 5).Value = '  -7.21%  '

$cell = $ws.Cells.Item(32, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.467'
$cell.Style = "Normal"
$ws.Cells.Item(32, 5).Value = '  -8.60%  '

$cell = $ws.Cells.Item(33, 4)
$cell.NumberFormat = "@"
$cell.Value = '3.460'
$cell.Style = "Normal"
$ws.Cells.Item(33, 5).Value = '  -8.05%  '

$cell = $ws.Cells.Item(34, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.646'
$cell.Style = "Normal"
$ws.Cells.Item(34, 5).Value = '  -4.48%  '

$cell = $ws.Cells.Item(35, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.004'
$cell.Style = "Normal"
$ws.Cells.Item(35, 5).Value = '  -5.26%  '

$ws.Cells.Item(36, 5).Value = '  -4.61%  '

$cell = $ws.Cells.Item(37, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.5857'
$cell.Style = "Normal"
$ws.Cells.Item(37, 5).Value = '  -8.90%  '

$cell = $ws.Cells.Item(38, 4)
$cell.NumberFormat = "@"
$cell.Value = '2.618'
$cell.Style = "Normal"
$ws.Cells.Item(38, 5).Value = '  -3.89%  '

$cell = $ws.Cells.Item(39, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.01582'
$cell.Style = "Normal"
$ws.Cells.Item(39, 5).Value = '  -7.41%  '

$cell = $ws.Cells.Item(40, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.070.05'
$cell.Style = "Normal"
$ws.Cells.Item(40, 5).Value = '  -6.46%  '

$cell = $ws.Cells.Item(41, 4)
$cell.NumberFormat = "@"
$cell.Value = '5.893'
$cell.Style = "Normal"
$ws.Cells.Item(41, 5).Value = '  -8.38%  '

$cell = $ws.Cells.Item(42, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.8490'
$cell.Style = "Normal"
$ws.Cells.Item(42, 5).Value = '  -2.20%  '

$ws.Cells.Item(43, 5).Value = '  -0.24%  '

$cell = $ws.Cells.Item(44, 4)
$cell.NumberFormat = "@"
$cell.Value = '98.03'
$cell.Style = "Normal"
$ws.Cells.Item(44, 5).Value = '  -2.76%  '

$cell = $ws.Cells.Item(45, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.843.24'
$cell.Style = "Normal"
$ws.Cells.Item(45, 5).Value = '  -6.11%  '

$cell = $ws.Cells.Item(46, 4)
$cell.NumberFormat = "@"
$cell.Value = '56.14'
$cell.Style = "Normal"
$ws.Cells.Item(46, 5).Value = '  -7.50%  '

$cell = $ws.Cells.Item(47, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.00000000105'
$cell.Style = "Normal"
$ws.Cells.Item(47, 5).Value = '  -5.99%  '

$cell = $ws.Cells.Item(48, 4)
$cell.NumberFormat = "@"
$cell.Value = '1.004'
$cell.Style = "Normal"
$ws.Cells.Item(48, 5).Value = '  -0.36%  '

$cell = $ws.Cells.Item(49, 4)
$cell.NumberFormat = "@"
$cell.Value = '8.066'
$cell.Style = "Normal"
$ws.Cells.Item(49, 5).Value = '  -3.93%  '

$cell = $ws.Cells.Item(50, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.4329'
$cell.Style = "Normal"
$ws.Cells.Item(50, 5).Value = '  -3.65%  '

$cell = $ws.Cells.Item(51, 4)
$cell.NumberFormat = "@"
$cell.Value = '0.05202'
$cell.Style = "Normal"
$ws.Cells.Item(51, 5).Value = '  -4.79%  '

